$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 6.449754000000001
$ws.Cells.Item(2, 8).Value = 19.349262
$ws.Cells.Item(2, 9).Value = 0.03479900749229446
$ws.Cells.Item(2, 10).Value = 0.03479900749229446
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.02270466666666667
$ws.Cells.Item(2, 14).Value = 0.06811400000000001
$ws.Cells.Item(2, 15).Value = 0.002206225855740089
$ws.Cells.Item(2, 16).Value = 0.002206225855740089
$ws.Cells.Item(2, 17).Value = 0.146439514652
$ws.Cells.Item(2, 18).Value = 1.317955631868
$ws.Cells.Item(2, 19).Value = 0.000076774470083593098304088065
$ws.Cells.Item(2, 20).Value = 0.000076774470083593098304088065

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 6.449754000000001
$ws.Cells.Item(3, 8).Value = 19.349262
$ws.Cells.Item(3, 9).Value = 0.03479900749229446
$ws.Cells.Item(3, 10).Value = 0.03479900749229446
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.02347533333333333
$ws.Cells.Item(3, 14).Value = 0.070426
$ws.Cells.Item(3, 15).Value = 0.002281111990432972
$ws.Cells.Item(3, 16).Value = 0.002281111990432972
$ws.Cells.Item(3, 17).Value = 0.151410125068
$ws.Cells.Item(3, 18).Value = 1.362691125612
$ws.Cells.Item(3, 19).Value = 0.000079380433245839732386665843
$ws.Cells.Item(3, 20).Value = 0.000079380433245839732386665843

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 6.449754000000001
$ws.Cells.Item(4, 8).Value = 19.349262
$ws.Cells.Item(4, 9).Value = 0.03479900749229446
$ws.Cells.Item(4, 10).Value = 0.03479900749229446
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 10.24499966666667
$ws.Cells.Item(4, 14).Value = 30.734999
$ws.Cells.Item(4, 15).Value = 0.9955126621538269
$ws.Cells.Item(4, 16).Value = 0.9955126621538269
$ws.Cells.Item(4, 17).Value = 66.07772758008201
$ws.Cells.Item(4, 18).Value = 594.6995482207382
$ws.Cells.Item(4, 19).Value = 0.03464285258896503
$ws.Cells.Item(4, 20).Value = 0.03464285258896503

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 123.027733
$ws.Cells.Item(5, 8).Value = 369.083199
$ws.Cells.Item(5, 9).Value = 0.663783921437469
$ws.Cells.Item(5, 10).Value = 0.6637839214374691
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.02270466666666667
$ws.Cells.Item(5, 14).Value = 0.06811400000000001
$ws.Cells.Item(5, 15).Value = 0.002206225855740089
$ws.Cells.Item(5, 16).Value = 0.002206225855740089
$ws.Cells.Item(5, 17).Value = 2.793303668520667
$ws.Cells.Item(5, 18).Value = 25.139733016686
$ws.Cells.Item(5, 19).Value = 0.001464457250099892
$ws.Cells.Item(5, 20).Value = 0.001464457250099892

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 123.027733
$ws.Cells.Item(6, 8).Value = 369.083199
$ws.Cells.Item(6, 9).Value = 0.663783921437469
$ws.Cells.Item(6, 10).Value = 0.6637839214374691
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.02347533333333333
$ws.Cells.Item(6, 14).Value = 0.070426
$ws.Cells.Item(6, 15).Value = 0.002281111990432972
$ws.Cells.Item(6, 16).Value = 0.002281111990432972
$ws.Cells.Item(6, 17).Value = 2.888117041419334
$ws.Cells.Item(6, 18).Value = 25.993053372774
$ws.Cells.Item(6, 19).Value = 0.001514165462247629
$ws.Cells.Item(6, 20).Value = 0.001514165462247629

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 123.027733
$ws.Cells.Item(7, 8).Value = 369.083199
$ws.Cells.Item(7, 9).Value = 0.663783921437469
$ws.Cells.Item(7, 10).Value = 0.6637839214374691
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 10.24499966666667
$ws.Cells.Item(7, 14).Value = 30.734999
$ws.Cells.Item(7, 15).Value = 0.9955126621538269
$ws.Cells.Item(7, 16).Value = 0.9955126621538269
$ws.Cells.Item(7, 17).Value = 1260.419083575756
$ws.Cells.Item(7, 18).Value = 11343.7717521818
$ws.Cells.Item(7, 19).Value = 0.6608052987251215
$ws.Cells.Item(7, 20).Value = 0.6608052987251216

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 55.79038633333334
$ws.Cells.Item(8, 8).Value = 167.371159
$ws.Cells.Item(8, 9).Value = 0.3010114916028843
$ws.Cells.Item(8, 10).Value = 0.3010114916028843
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.02270466666666667
$ws.Cells.Item(8, 14).Value = 0.06811400000000001
$ws.Cells.Item(8, 15).Value = 0.002206225855740089
$ws.Cells.Item(8, 16).Value = 0.002206225855740089
$ws.Cells.Item(8, 17).Value = 1.266702124902889
$ws.Cells.Item(8, 18).Value = 11.400319124126
$ws.Cells.Item(8, 19).Value = 0.0006640993356491738
$ws.Cells.Item(8, 20).Value = 0.0006640993356491738

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 55.79038633333334
$ws.Cells.Item(9, 8).Value = 167.371159
$ws.Cells.Item(9, 9).Value = 0.3010114916028843
$ws.Cells.Item(9, 10).Value = 0.3010114916028843
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.02347533333333333
$ws.Cells.Item(9, 14).Value = 0.070426
$ws.Cells.Item(9, 15).Value = 0.002281111990432972
$ws.Cells.Item(9, 16).Value = 0.002281111990432972
$ws.Cells.Item(9, 17).Value = 1.309697915970444
$ws.Cells.Item(9, 18).Value = 11.787281243734
$ws.Cells.Item(9, 19).Value = 0.0006866409227534533
$ws.Cells.Item(9, 20).Value = 0.0006866409227534533

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 55.79038633333334
$ws.Cells.Item(10, 8).Value = 167.371159
$ws.Cells.Item(10, 9).Value = 0.3010114916028843
$ws.Cells.Item(10, 10).Value = 0.3010114916028843
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 10.24499966666667
$ws.Cells.Item(10, 14).Value = 30.734999
$ws.Cells.Item(10, 15).Value = 0.9955126621538269
$ws.Cells.Item(10, 16).Value = 0.9955126621538269
$ws.Cells.Item(10, 17).Value = 571.5724893882046
$ws.Cells.Item(10, 18).Value = 5144.152404493841
$ws.Cells.Item(10, 19).Value = 0.2996607513444816
$ws.Cells.Item(10, 20).Value = 0.2996607513444816

$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.07517133333333333
$ws.Cells.Item(11, 8).Value = 0.225514
$ws.Cells.Item(11, 9).Value = 0.0004055794673521549
$ws.Cells.Item(11, 10).Value = 0.000405579467352155
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.02270466666666667
$ws.Cells.Item(11, 14).Value = 0.06811400000000001
$ws.Cells.Item(11, 15).Value = 0.002206225855740089
$ws.Cells.Item(11, 16).Value = 0.002206225855740089
$ws.Cells.Item(11, 17).Value = 0.001706740066222222
$ws.Cells.Item(11, 18).Value = 0.015360660596
$ws.Cells.Item(11, 19).Value = 0.000000894799907429617306708685
$ws.Cells.Item(11, 20).Value = 0.000000894799907429617412587804

$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.07517133333333333
$ws.Cells.Item(12, 8).Value = 0.225514
$ws.Cells.Item(12, 9).Value = 0.0004055794673521549
$ws.Cells.Item(12, 10).Value = 0.000405579467352155
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 0.6666666666666666
$ws.Cells.Item(12, 13).Value = 0.02347533333333333
$ws.Cells.Item(12, 14).Value = 0.070426
$ws.Cells.Item(12, 15).Value = 0.002281111990432972
$ws.Cells.Item(12, 16).Value = 0.002281111990432972
$ws.Cells.Item(12, 17).Value = 0.001764672107111111
$ws.Cells.Item(12, 18).Value = 0.015882048964
$ws.Cells.Item(12, 19).Value = 0.000000925172186050418743539733
$ws.Cells.Item(12, 20).Value = 0.00000092517218605041895529797

$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.07517133333333333
$ws.Cells.Item(13, 8).Value = 0.225514
$ws.Cells.Item(13, 9).Value = 0.0004055794673521549
$ws.Cells.Item(13, 10).Value = 0.000405579467352155
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 10.24499966666667
$ws.Cells.Item(13, 14).Value = 30.734999
$ws.Cells.Item(13, 15).Value = 0.9955126621538269
$ws.Cells.Item(13, 16).Value = 0.9955126621538269
$ws.Cells.Item(13, 17).Value = 0.7701302849428888
$ws.Cells.Item(13, 18).Value = 6.931172564486
$ws.Cells.Item(13, 19).Value = 0.0004037594952586748
$ws.Cells.Item(13, 20).Value = 0.0004037594952586749
